$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1528.9
$ws.Range("J70").Value = 1612.7142
$ws.Range("L70").Value = 4838.142599999999
$ws.Range("N70").Value = -5378.142599999999
$ws.Range("H73").Value = 1528.9
$ws.Range("J73").Value = 1612.7142
$ws.Range("L73").Value = 4838.142599999999
$ws.Range("N73").Value = -6710.142599999999
$ws.Range("H112").Value = 3403.3333
$ws.Range("I112").Value = 2550
$ws.Range("J112").Value = 3453.5293
$ws.Range("K112").Value = 7650
$ws.Range("L112").Value = 10360.5879
$ws.Range("M112").Value = -6542
$ws.Range("N112").Value = -12576.5879
$ws.Range("H116").Value = 2900.4546
$ws.Range("I116").Value = 2822.7778
$ws.Range("J116").Value = 3250
$ws.Range("K116").Value = 2822.7778
$ws.Range("L116").Value = 3250
$ws.Range("M116").Value = 619.2222000000002
$ws.Range("N116").Value = -10134
$ws.Range("H132").Value = 5351.3105
$ws.Range("I132").Value = 5073.174
$ws.Range("K132").Value = 15219.522
$ws.Range("M132").Value = -12689.522
$ws.Range("H137").Value = 1086.3636
$ws.Range("I137").Value = 476.81818
$ws.Range("J137").Value = 1695.909
$ws.Range("K137").Value = 1430.45454
$ws.Range("L137").Value = 5087.727000000001
$ws.Range("M137").Value = 1119.54546
$ws.Range("N137").Value = -10187.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13300
$ws.Range("I2").Value = 14996.857
$ws.Range("J2").Value = 1422
$ws.Range("K2").Value = 14996.857
$ws.Range("L2").Value = 1422
$ws.Range("M2").Value = -14883.857
$ws.Range("N2").Value = -1648
$ws.Range("H32").Value = 667379.4399999999
$ws.Range("I32").Value = 844345.1
$ws.Range("K32").Value = 844345.1
$ws.Range("M32").Value = -844058.1
$ws.Range("H110").Value = 2954.2273
$ws.Range("I110").Value = 2704.2104
$ws.Range("J110").Value = 4537.6665
$ws.Range("K110").Value = 2704.2104
$ws.Range("L110").Value = 4537.6665
$ws.Range("M110").Value = -659.2103999999999
$ws.Range("N110").Value = -8627.666499999999
$ws.Range("H116").Value = 13300
$ws.Range("I116").Value = 14996.857
$ws.Range("J116").Value = 1422
$ws.Range("K116").Value = 14996.857
$ws.Range("L116").Value = 1422
$ws.Range("M116").Value = -12702.857
$ws.Range("N116").Value = -6010
$ws.Range("H122").Value = 1192.35
$ws.Range("I122").Value = 1071.6875
$ws.Range("J122").Value = 1675
$ws.Range("K122").Value = 3215.0625
$ws.Range("L122").Value = 5025
$ws.Range("M122").Value = -765.0625
$ws.Range("N122").Value = -9925

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13300
$ws.Range("I3").Value = 14996.857
$ws.Range("J3").Value = 1422
$ws.Range("K3").Value = 14996.857
$ws.Range("L3").Value = 1422
$ws.Range("M3").Value = -14882.857
$ws.Range("N3").Value = -1650
$ws.Range("H20").Value = 2144.75
$ws.Range("I20").Value = 2000.8572
$ws.Range("J20").Value = 2256.6667
$ws.Range("K20").Value = 2000.8572
$ws.Range("L20").Value = 2256.6667
$ws.Range("M20").Value = -1753.8572
$ws.Range("N20").Value = -2750.6667
$ws.Range("I94").Value = 918.2857
$ws.Range("J94").Value = 819
$ws.Range("K94").Value = 918.2857
$ws.Range("L94").Value = 819
$ws.Range("M94").Value = -467.2857
$ws.Range("N94").Value = -1721
$ws.Range("H134").Value = 2645.122
$ws.Range("I134").Value = 2608.4827
$ws.Range("K134").Value = 7825.4481
$ws.Range("M134").Value = -5290.4481

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 3950
$ws.Range("I29").Value = 3900
$ws.Range("K29").Value = 3900
$ws.Range("M29").Value = -3607
$ws.Range("H31").Value = 1424.7273
$ws.Range("I31").Value = 1298.9166
$ws.Range("J31").Value = 1575.7
$ws.Range("K31").Value = 1298.9166
$ws.Range("L31").Value = 1575.7
$ws.Range("M31").Value = -1003.9166
$ws.Range("N31").Value = -2165.7
$ws.Range("H34").Value = 1424.7273
$ws.Range("I34").Value = 1298.9166
$ws.Range("J34").Value = 1575.7
$ws.Range("K34").Value = 1298.9166
$ws.Range("L34").Value = 1575.7
$ws.Range("M34").Value = -1096.9166
$ws.Range("N34").Value = -1979.7
$ws.Range("H35").Value = 1912.5
$ws.Range("I35").Value = 1912.5
$ws.Range("K35").Value = 1912.5
$ws.Range("M35").Value = -1618.5
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 909.375
$ws.Range("I38").Value = 100
$ws.Range("J38").Value = 963.3333
$ws.Range("K38").Value = 300
$ws.Range("L38").Value = 2889.9999
$ws.Range("M38").Value = 47
$ws.Range("N38").Value = -3583.9999
$ws.Range("H68").Value = 908.6415
$ws.Range("I68").Value = 659
$ws.Range("J68").Value = 966.6977000000001
$ws.Range("K68").Value = 1977
$ws.Range("L68").Value = 2900.0931
$ws.Range("M68").Value = -1166
$ws.Range("N68").Value = -4522.0931
$ws.Range("H71").Value = 908.6415
$ws.Range("I71").Value = 659
$ws.Range("J71").Value = 966.6977000000001
$ws.Range("K71").Value = 5931
$ws.Range("L71").Value = 8700.2793
$ws.Range("M71").Value = -1875
$ws.Range("N71").Value = -16812.2793
$ws.Range("H113").Value = 1570
$ws.Range("I113").Value = 750
$ws.Range("J113").Value = 1775
$ws.Range("K113").Value = 2250
$ws.Range("L113").Value = 5325
$ws.Range("M113").Value = -80
$ws.Range("N113").Value = -9665
$ws.Range("H118").Value = 2835
$ws.Range("J118").Value = 2835
$ws.Range("L118").Value = 8505
$ws.Range("N118").Value = -10991
$ws.Range("H132").Value = 1196.15
$ws.Range("I132").Value = 1338
$ws.Range("J132").Value = 1101.5834
$ws.Range("K132").Value = 12042
$ws.Range("L132").Value = 9914.250599999999
$ws.Range("M132").Value = -9512
$ws.Range("N132").Value = -14974.2506

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3452.0557
$ws.Range("I132").Value = 4203.4287
$ws.Range("J132").Value = 2973.9092
$ws.Range("K132").Value = 12610.2861
$ws.Range("L132").Value = 8921.7276
$ws.Range("M132").Value = -10080.2861
$ws.Range("N132").Value = -13981.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2689.4546
$ws.Range("I7").Value = 2400.6667
$ws.Range("J7").Value = 3036
$ws.Range("K7").Value = 2400.6667
$ws.Range("L7").Value = 3036
$ws.Range("M7").Value = -2288.6667
$ws.Range("N7").Value = -3260
$ws.Range("H61").Value = 4336.3335
$ws.Range("I61").Value = 4398.1665
$ws.Range("J61").Value = 4253.8887
$ws.Range("K61").Value = 4398.1665
$ws.Range("L61").Value = 4253.8887
$ws.Range("M61").Value = -4196.1665
$ws.Range("N61").Value = -4657.8887
$ws.Range("H113").Value = 4336.3335
$ws.Range("I113").Value = 4398.1665
$ws.Range("J113").Value = 4253.8887
$ws.Range("K113").Value = 4398.1665
$ws.Range("L113").Value = 4253.8887
$ws.Range("M113").Value = -2228.1665
$ws.Range("N113").Value = -8593.8887
$ws.Range("H126").Value = 2689.4546
$ws.Range("I126").Value = 2400.6667
$ws.Range("J126").Value = 3036
$ws.Range("K126").Value = 7202.000100000001
$ws.Range("L126").Value = 9108
$ws.Range("M126").Value = -4732.000100000001
$ws.Range("N126").Value = -14048
$ws.Range("H132").Value = 3638.3784
$ws.Range("I132").Value = 3218.348
$ws.Range("K132").Value = 9655.044
$ws.Range("M132").Value = -7125.044

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1474.2142
$ws.Range("I122").Value = 1356.4615
$ws.Range("K122").Value = 4069.3845
$ws.Range("M122").Value = -1619.3845
$ws.Range("H126").Value = 1369
$ws.Range("I126").Value = 1147.2222
$ws.Range("J126").Value = 1701.6666
$ws.Range("K126").Value = 3441.6666
$ws.Range("L126").Value = 5104.9998
$ws.Range("M126").Value = -971.6665999999996
$ws.Range("N126").Value = -10044.9998
